# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# The "Date" column ("5-5-2013-14") was off by a day - correct it to "2014-05-05"
# for every data row on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow  = $usedRange.Row
$lastRow   = $usedRange.Rows.Count + $firstRow - 1
$firstCol  = $usedRange.Column
$lastCol   = $usedRange.Columns.Count + $firstCol - 1

# Locate the "Date" column from the header row instead of hard-coding it.
$dateCol = 0
for ($c = $firstCol; $c -le $lastCol; $c++) {
    if ($ws.Cells.Item($firstRow, $c).Text -eq "Date") {
        $dateCol = $c
        break
    }
}

if ($dateCol -gt 0) {
    $oldValue = "5-5-2013-14"
    $newValue = "2014-05-05"

    # Force the column to be treated as text first so Excel doesn't reinterpret
    # the corrected value (which looks like an ISO date) as a date serial number;
    # we want the literal string "2014-05-05" kept, just like the original text.
    $dateRange = $ws.Range($ws.Cells.Item($firstRow + 1, $dateCol), $ws.Cells.Item($lastRow, $dateCol))
    $dateRange.NumberFormat = "@"

    for ($r = $firstRow + 1; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, $dateCol)
        if ($cell.Text -eq $oldValue) {
            $cell.Value = $newValue
        }
    }
}
